$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B26").Value = 6489
$ws.Range("C26").Value = 1010
$ws.Range("D26").Value = 6046302
$ws.Range("E26").Value = 931.7771613499768
$ws.Range("F26").Value = 9.667061010647293
$ws.Range("G26").Value = 7.218683651804669
$ws.Range("H26").Value = 25.91444284808993
